$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix cell values for A3, A4, A5 to all use "example@domain.com"
$ws.Range("A3").Value2 = "example@domain.com"
$ws.Range("A4").Value2 = "example@domain.com"
$ws.Range("A5").Value2 = "example@domain.com"

# Remove the old hyperlinks and rebuild: A2 alone, A3:A11 as one merged range
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:example@domain.com")
$ws.Hyperlinks.Add($ws.Range("A3:A11"), "mailto:example@domain.com", "", "", "example@domain.com")
$ws.Range("A2:A11").Style = "Hyperlink"

# Delete rows 12 and 13 (example2@domain.com / 1619DKV@ARMORSEC.XYZ rows)
$ws.Rows("12:13").Delete()

# Update selection to match target (A2:A11 selected, active cell A2)
$ws.Range("A2:A11").Select()

Write-Host "done"
